$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "legal_basis" translation row (new row 34, before the
#     "interstitial-block" row) ---
$ws.Range("A34").EntireRow.Insert()
$ws.Range("A34").Value2 = "legal_basis"
$ws.Range("B34").Value2 = "Notice is given under section 1 of the Statistics of Trade Act 1947."
$ws.Range("C34").Value2 = "Rhoddir rhybudd o dan adran 1 o Ddeddf Ystadegau Masnach 1947."
$ws.Range("C34").Font.Color = 2171169

# --- Insert "feeling-bad-answer" translation row (new row 41, before the
#     "feeling-answer [question guidance]" row) ---
$ws.Range("A41").EntireRow.Insert()
$ws.Range("A41").Value2 = "feeling-bad-answer"
$ws.Range("B41").Value2 = "Specify why answering for yourself is bad"
$ws.Range("C41").Value2 = "Nodwch pam mae ateb drosti eich hun yn wael"
$ws.Range("C41").Font.ThemeColor = 1

# --- The two row inserts above push the existing GDP hyperlink (originally
#     anchored at C37) down to C38; rebuild it there. ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C38"), "https://www.ons.gov.uk/economy/grossdomesticproductgdp")
